$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CCSaleData")
$ws.Range("A2").NumberFormat = "General"
